$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new product rows (6-11), columns A-G ---
$ws.Range("A6").Value = "'Melão"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "'100"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'14/05/2025"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'super"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'Alimento"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'2"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'boa`n"
$ws.Range("G6").Style = "Normal"
$ws.Range("A7").Value = "'Melão "
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "'100"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'03/06/2025"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'super"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'Alimento"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'2"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'é isso`n"
$ws.Range("G7").Style = "Normal"
$ws.Range("A8").Value = "'Melão"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "'100"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'05/06/2025"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'Super"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'Alimento"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'2"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'é isso`n"
$ws.Range("G8").Style = "Normal"
$ws.Range("A9").Value = "'Melão"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = "'100"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'05/06/2025"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'Super"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'Alimento"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'2"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = "'É isso`n"
$ws.Range("G9").Style = "Normal"
$ws.Range("A10").Value = "'Melão"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = "'100"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'06/06/2025"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'Super"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'Alimento"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'2"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = "'É isso`n"
$ws.Range("G10").Style = "Normal"
$ws.Range("A11").Value = "'Melão"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = "'100"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'06/06/2025"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'Super"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'Alimento"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "'2"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = "'É isso`n"
$ws.Range("G11").Style = "Normal"

# --- Avoid auto row-height changes triggered by embedded newlines ---
foreach ($r in 6..11) {
  $ws.Rows.Item($r).AutoFit()
}

# --- Extend the table (ListObject) range to cover the new data ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G11"))

Write-Output "Edit complete"
